$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" footer date field from
#    03/09/2020 -> 30/11/2020 everywhere it appears: the slide master
#    and every custom (slide) layout's Date Placeholder.
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "03/09/2020") {
                $tr.Text = "30/11/2020"
            }
        }
    }
}

$sm = $p.SlideMaster
Update-DatePlaceholders $sm.Shapes

for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $layout = $sm.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Re-style the legend text on slide 1 (inside "Group 1"):
#    - "Type of citation:" drops from 11pt to 10pt
#    - every legend label switches from Arial to Times New Roman
#      (both the Latin and Complex-Script typeface)
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
$group = $s.Shapes.Item(1)

for ($i = 1; $i -le $group.GroupItems.Count; $i++) {
    $sh = $group.GroupItems.Item($i)
    if (-not $sh.HasTextFrame) { continue }

    $tr = $sh.TextFrame.TextRange
    if ($tr.Text.Length -eq 0) { continue }

    if ($tr.Text -eq "Type of citation:") {
        $tr.Font.Size = 10
    }

    $tr.Font.Name = "Times New Roman"
    $tr.Font.NameComplexScript = "Times New Roman"
}
